# Daily attendance processing - 2026-01-12 09:51:13
# Swap "System, <user>" to "<user>, System" in the "Recorded By" column (G)
# for every data row, leaving already-correct / multi-party / System-only
# values untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    if ($val -match "^System, ([^,]+)$") {
        $other = $matches[1]
        if (-not $other.Contains("backdoor")) {
            $cell.Value = $other + ", System"
        }
    }
}

Write-Host "Recorded By column normalized"
